# Update NATMI Bmp2-Acvr2b LR-pair TPM-derived metrics (YoungD7) to reflect new TPM computation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.729584333333333
$ws.Range("H2").Value = 5.188753
$ws.Range("I2").Value = 0.2476387648475193
$ws.Range("J2").Value = 0.2476387648475193
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.110264333333333
$ws.Range("N2").Value = 6.330793
$ws.Range("O2").Value = 0.3832041185227171
$ws.Range("P2").Value = 0.3832041185227171
$ws.Range("Q2").Value = 3.649880130125445
$ws.Range("R2").Value = 32.848921171129
$ws.Range("S2").Value = 0.09489619459544804
$ws.Range("T2").Value = 0.09489619459544804

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.729584333333333
$ws.Range("H3").Value = 5.188753
$ws.Range("I3").Value = 0.2476387648475193
$ws.Range("J3").Value = 0.2476387648475193
$ws.Range("O3").Value = 0.0946183755984393
$ws.Range("P3").Value = 0.0946183755984393
$ws.Range("Q3").Value = 0.9012056821644443
$ws.Range("R3").Value = 8.110851139479999
$ws.Range("S3").Value = 0.02343117766507617
$ws.Range("T3").Value = 0.02343117766507617

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.729584333333333
$ws.Range("H4").Value = 5.188753
$ws.Range("I4").Value = 0.2476387648475193
$ws.Range("J4").Value = 0.2476387648475193
$ws.Range("M4").Value = 2.065388333333333
$ws.Range("N4").Value = 6.196165
$ws.Range("O4").Value = 0.3750550597762889
$ws.Range("P4").Value = 0.3750550597762889
$ws.Range("Q4").Value = 3.572263303582778
$ws.Range("R4").Value = 32.150369732245
$ws.Range("S4").Value = 0.0928781717528127
$ws.Range("T4").Value = 0.0928781717528127

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.729584333333333
$ws.Range("H5").Value = 5.188753
$ws.Range("I5").Value = 0.2476387648475193
$ws.Range("J5").Value = 0.2476387648475193
$ws.Range("M5").Value = 0.8101876666666666
$ws.Range("N5").Value = 2.430563
$ws.Range("O5").Value = 0.1471224461025547
$ws.Range("P5").Value = 0.1471224461025547
$ws.Range("Q5").Value = 1.401287895326556
$ws.Range("R5").Value = 12.611591057939
$ws.Range("S5").Value = 0.03643322083418239
$ws.Range("T5").Value = 0.03643322083418239

# Row 6
$ws.Range("I6").Value = 0.2307941364328804
$ws.Range("J6").Value = 0.2307941364328804
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.110264333333333
$ws.Range("N6").Value = 6.330793
$ws.Range("O6").Value = 0.3832041185227171
$ws.Range("P6").Value = 0.3832041185227171
$ws.Range("Q6").Value = 3.401611751837445
$ws.Range("R6").Value = 30.614505766537
$ws.Range("S6").Value = 0.08844126361197363
$ws.Range("T6").Value = 0.08844126361197363

# Row 7
$ws.Range("I7").Value = 0.2307941364328804
$ws.Range("J7").Value = 0.2307941364328804
$ws.Range("O7").Value = 0.0946183755984393
$ws.Range("P7").Value = 0.0946183755984393
$ws.Range("S7").Value = 0.02183736628692372
$ws.Range("T7").Value = 0.02183736628692372

# Row 8
$ws.Range("I8").Value = 0.2307941364328804
$ws.Range("J8").Value = 0.2307941364328804
$ws.Range("M8").Value = 2.065388333333333
$ws.Range("N8").Value = 6.196165
$ws.Range("O8").Value = 0.3750550597762889
$ws.Range("P8").Value = 0.3750550597762889
$ws.Range("Q8").Value = 3.329274496942778
$ws.Range("R8").Value = 29.963470472485
$ws.Range("S8").Value = 0.08656050863585092
$ws.Range("T8").Value = 0.08656050863585092

# Row 9
$ws.Range("I9").Value = 0.2307941364328804
$ws.Range("J9").Value = 0.2307941364328804
$ws.Range("M9").Value = 0.8101876666666666
$ws.Range("N9").Value = 2.430563
$ws.Range("O9").Value = 0.1471224461025547
$ws.Range("P9").Value = 0.1471224461025547
$ws.Range("Q9").Value = 1.305970936718556
$ws.Range("R9").Value = 11.753738430467
$ws.Range("S9").Value = 0.0339549978981321
$ws.Range("T9").Value = 0.0339549978981321

# Row 10
$ws.Range("G10").Value = 2.743651333333334
$ws.Range("H10").Value = 8.230954000000001
$ws.Range("I10").Value = 0.3928310486309039
$ws.Range("J10").Value = 0.3928310486309038
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.110264333333333
$ws.Range("N10").Value = 6.330793
$ws.Range("O10").Value = 0.3832041185227171
$ws.Range("P10").Value = 0.3832041185227171
$ws.Range("Q10").Value = 5.789829551835778
$ws.Range("R10").Value = 52.108465966522
$ws.Range("S10").Value = 0.1505344757189601
$ws.Range("T10").Value = 0.1505344757189601

# Row 11
$ws.Range("G11").Value = 2.743651333333334
$ws.Range("H11").Value = 8.230954000000001
$ws.Range("I11").Value = 0.3928310486309039
$ws.Range("J11").Value = 0.3928310486309038
$ws.Range("O11").Value = 0.0946183755984393
$ws.Range("P11").Value = 0.0946183755984393
$ws.Range("Q11").Value = 1.429588672737778
$ws.Range("R11").Value = 12.86629805464
$ws.Range("S11").Value = 0.03716903570608764
$ws.Range("T11").Value = 0.03716903570608764

# Row 12
$ws.Range("G12").Value = 2.743651333333334
$ws.Range("H12").Value = 8.230954000000001
$ws.Range("I12").Value = 0.3928310486309039
$ws.Range("J12").Value = 0.3928310486309038
$ws.Range("M12").Value = 2.065388333333333
$ws.Range("N12").Value = 6.196165
$ws.Range("O12").Value = 0.3750550597762889
$ws.Range("P12").Value = 0.3750550597762889
$ws.Range("Q12").Value = 5.666705454601112
$ws.Range("R12").Value = 51.00034909141
$ws.Range("S12").Value = 0.1473332724262459
$ws.Range("T12").Value = 0.1473332724262459

# Row 13
$ws.Range("G13").Value = 2.743651333333334
$ws.Range("H13").Value = 8.230954000000001
$ws.Range("I13").Value = 0.3928310486309039
$ws.Range("J13").Value = 0.3928310486309038
$ws.Range("M13").Value = 0.8101876666666666
$ws.Range("N13").Value = 2.430563
$ws.Range("O13").Value = 0.1471224461025547
$ws.Range("P13").Value = 0.1471224461025547
$ws.Range("Q13").Value = 2.222872471900223
$ws.Range("R13").Value = 20.005852247102
$ws.Range("S13").Value = 0.05779426477961021
$ws.Range("T13").Value = 0.05779426477961021

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.8991316666666668
$ws.Range("H14").Value = 2.697395
$ws.Range("I14").Value = 0.1287360500886965
$ws.Range("J14").Value = 0.1287360500886965
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.110264333333333
$ws.Range("N14").Value = 6.330793
$ws.Range("O14").Value = 0.3832041185227171
$ws.Range("P14").Value = 0.3832041185227171
$ws.Range("Q14").Value = 1.897405487137222
$ws.Range("R14").Value = 17.076649384235
$ws.Range("S14").Value = 0.0493321845963353
$ws.Range("T14").Value = 0.0493321845963353

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.8991316666666668
$ws.Range("H15").Value = 2.697395
$ws.Range("I15").Value = 0.1287360500886965
$ws.Range("J15").Value = 0.1287360500886965
$ws.Range("O15").Value = 0.0946183755984393
$ws.Range("P15").Value = 0.0946183755984393
$ws.Range("Q15").Value = 0.4684955520222222
$ws.Range("R15").Value = 4.2164599682
$ws.Range("S15").Value = 0.01218079594035178
$ws.Range("T15").Value = 0.01218079594035178

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.8991316666666668
$ws.Range("H16").Value = 2.697395
$ws.Range("I16").Value = 0.1287360500886965
$ws.Range("J16").Value = 0.1287360500886965
$ws.Range("M16").Value = 2.065388333333333
$ws.Range("N16").Value = 6.196165
$ws.Range("O16").Value = 0.3750550597762889
$ws.Range("P16").Value = 0.3750550597762889
$ws.Range("Q16").Value = 1.857056054463889
$ws.Range("R16").Value = 16.713504490175
$ws.Range("S16").Value = 0.04828310696137939
$ws.Range("T16").Value = 0.04828310696137939

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.8991316666666668
$ws.Range("H17").Value = 2.697395
$ws.Range("I17").Value = 0.1287360500886965
$ws.Range("J17").Value = 0.1287360500886965
$ws.Range("M17").Value = 0.8101876666666666
$ws.Range("N17").Value = 2.430563
$ws.Range("O17").Value = 0.1471224461025547
$ws.Range("P17").Value = 0.1471224461025547
$ws.Range("Q17").Value = 0.7284653870427779
$ws.Range("R17").Value = 6.556188483385
$ws.Range("S17").Value = 0.01893996259063004
$ws.Range("T17").Value = 0.01893996259063004
